$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.542.34"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.081.15"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.635"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.391"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "2.388.91"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "2.079.82"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "37.499.95"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  +6.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0232"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.54%  "
$ws.Range("D46").Value = "1.459.72"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "2.272.73"
$ws.Range("E51").Value = "  +0.19%  "
